$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2173913043478261
$ws.Range("C2").Value = 0.5072463768115942
$ws.Range("J2").Value = 0.02608695652173913
$ws.Range("P2").Value = 0.1449275362318841
$ws.Range("S2").Value = 0.1043478260869565
$ws.Range("B3").Value = 0.005494505494505495
$ws.Range("C3").Value = 0.03296703296703297
$ws.Range("J3").Value = 0.02197802197802198
$ws.Range("P3").Value = 0.7142857142857143
$ws.Range("S3").Value = 0.2252747252747253
$ws.Range("J4").Value = 0.04166666666666666
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.2916666666666667
$ws.Range("B6").Value = 0.05058365758754864
$ws.Range("D6").Value = 0.007782101167315175
$ws.Range("E6").Value = 0.003891050583657588
$ws.Range("F6").Value = 0.06614785992217899
$ws.Range("J6").Value = 0.245136186770428
$ws.Range("O6").Value = 0.01167315175097276
$ws.Range("Q6").Value = 0.1556420233463035
$ws.Range("R6").Value = 0.07782101167315175
$ws.Range("S6").Value = 0.3813229571984436
$ws.Range("B7").Value = 0.08056872037914692
$ws.Range("D7").Value = 0.02369668246445497
$ws.Range("E7").Value = 0.004739336492890996
$ws.Range("F7").Value = 0.06161137440758294
$ws.Range("J7").Value = 0.1990521327014218
$ws.Range("O7").Value = 0.01895734597156398
$ws.Range("Q7").Value = 0.1090047393364929
$ws.Range("R7").Value = 0.1232227488151659
$ws.Range("S7").Value = 0.3791469194312796
$ws.Range("B8").Value = 0.08058608058608059
$ws.Range("D8").Value = 0.003663003663003663
$ws.Range("F8").Value = 0.05677655677655678
$ws.Range("J8").Value = 0.1007326007326007
$ws.Range("O8").Value = 0.007326007326007326
$ws.Range("Q8").Value = 0.1794871794871795
$ws.Range("R8").Value = 0.119047619047619
$ws.Range("S8").Value = 0.4523809523809524
$ws.Range("B9").Value = 0.1106194690265487
$ws.Range("D9").Value = 0.01769911504424779
$ws.Range("F9").Value = 0.06194690265486726
$ws.Range("J9").Value = 0.1150442477876106
$ws.Range("O9").Value = 0.004424778761061947
$ws.Range("Q9").Value = 0.1548672566371681
$ws.Range("R9").Value = 0.06637168141592921
$ws.Range("S9").Value = 0.4690265486725664
$ws.Range("B10").Value = 0.1091525423728814
$ws.Range("D10").Value = 0.02576271186440678
$ws.Range("F10").Value = 0.07186440677966102
$ws.Range("J10").Value = 0.1349152542372881
$ws.Range("O10").Value = 0.008813559322033898
$ws.Range("Q10").Value = 0.1857627118644068
$ws.Range("R10").Value = 0.08949152542372882
$ws.Range("S10").Value = 0.3742372881355932
$ws.Range("G11").Value = 0.1780415430267062
$ws.Range("J11").Value = 0.0712166172106825
$ws.Range("K11").Value = 0.2344213649851632
$ws.Range("L11").Value = 0.5074183976261127
$ws.Range("S11").Value = 0.008902077151335312
$ws.Range("G12").Value = 0.7413793103448276
$ws.Range("J12").Value = 0.2298850574712644
$ws.Range("L12").Value = 0.005747126436781609
$ws.Range("S12").Value = 0.02298850574712644
$ws.Range("G13").Value = 0.5652173913043478
$ws.Range("J13").Value = 0.391304347826087
$ws.Range("S13").Value = 0.04347826086956522
$ws.Range("F15").Value = 0.00975609756097561
$ws.Range("H15").Value = 0.2048780487804878
$ws.Range("I15").Value = 0.08292682926829269
$ws.Range("J15").Value = 0.3707317073170732
$ws.Range("K15").Value = 0.04878048780487805
$ws.Range("M15").Value = 0.00975609756097561
$ws.Range("O15").Value = 0.09268292682926829
$ws.Range("S15").Value = 0.1804878048780488
$ws.Range("F16").Value = 0.004901960784313725
$ws.Range("H16").Value = 0.2107843137254902
$ws.Range("I16").Value = 0.08823529411764706
$ws.Range("J16").Value = 0.4019607843137255
$ws.Range("K16").Value = 0.08823529411764706
$ws.Range("M16").Value = 0.02450980392156863
$ws.Range("O16").Value = 0.0392156862745098
$ws.Range("S16").Value = 0.142156862745098
$ws.Range("F17").Value = 0.01502145922746781
$ws.Range("H17").Value = 0.2017167381974249
$ws.Range("I17").Value = 0.1072961373390558
$ws.Range("J17").Value = 0.4055793991416309
$ws.Range("K17").Value = 0.09012875536480687
$ws.Range("M17").Value = 0.01716738197424893
$ws.Range("N17").Value = 0.002145922746781116
$ws.Range("O17").Value = 0.04721030042918455
$ws.Range("S17").Value = 0.1137339055793991
$ws.Range("F18").Value = 0.01937984496124031
$ws.Range("H18").Value = 0.189922480620155
$ws.Range("I18").Value = 0.08914728682170543
$ws.Range("J18").Value = 0.4224806201550387
$ws.Range("K18").Value = 0.09689922480620156
$ws.Range("M18").Value = 0.02325581395348837
$ws.Range("N18").Value = 0.003875968992248062
$ws.Range("O18").Value = 0.06201550387596899
$ws.Range("S18").Value = 0.09302325581395349
$ws.Range("F19").Value = 0.01708817498291183
$ws.Range("H19").Value = 0.2228298017771702
$ws.Range("I19").Value = 0.07997265892002735
$ws.Range("J19").Value = 0.380041011619959
$ws.Range("K19").Value = 0.1093643198906357
$ws.Range("M19").Value = 0.01845522898154477
$ws.Range("N19").Value = 0.001367053998632946
$ws.Range("O19").Value = 0.06288448393711552
$ws.Range("S19").Value = 0.1079972658920027
